# Insert a new data row at row 10 (pushing existing rows 10..107 down to 11..108)
# and populate it with the new weekly price record for Acelga.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(10).Insert()

$ws.Cells.Item(10, 1).Value = 1
$ws.Cells.Item(10, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(10, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(10, 4).Value = 45230
$ws.Cells.Item(10, 5).Value = 15
$ws.Cells.Item(10, 6).Value = 100112009
$ws.Cells.Item(10, 7).Value = "Acelga"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 400
$ws.Cells.Item(10, 11).Value = 1000
$ws.Cells.Item(10, 12).Value = 1200
$ws.Cells.Item(10, 13).Value = 1125
$ws.Cells.Item(10, 14).Value = "`$/atado 2,5 a 3 kilos"
$ws.Cells.Item(10, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(10, 16).Value = 375
$ws.Cells.Item(10, 17).Value = 3
$ws.Cells.Item(10, 18).Value = "Hortaliza"
